$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.83"
$ws.Range("E2").Value = "'-0.62%"

$ws.Range("D3").Value = "'38.86"
$ws.Range("E3").Value = "'6.83%"

$ws.Range("D4").Value = "'5.111"
$ws.Range("E4").Value = "'1.04%"

$ws.Range("D5").Value = "'0.08070"
$ws.Range("E5").Value = "'-0.09%"

$ws.Range("D6").Value = "'1.931"
$ws.Range("E6").Value = "'-1.84%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.208"
$ws.Range("E7").Value = "'0.97%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.038"
$ws.Range("E8").Value = "'2.60%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9272"
$ws.Range("E9").Value = "'-0.07%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1429"
$ws.Range("E10").Value = "'-3.66%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1913"
$ws.Range("E11").Value = "'-1.35%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09043"
$ws.Range("E12").Value = "'-0.81%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03510"
$ws.Range("E13").Value = "'-0.57%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09760"
$ws.Range("E14").Value = "'-1.04%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001392"
$ws.Range("E15").Value = "'-1.19%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005875"
$ws.Range("E16").Value = "'-5.20%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.762"
$ws.Range("E17").Value = "'-2.05%"

$ws.Range("E18").Value = "'-1.99%"

$ws.Range("D19").Value = "'0.3461"
$ws.Range("E19").Value = "'0.38%"

$ws.Range("E20").Value = "'-0.01%"

$ws.Range("D21").Value = "'4.688"
$ws.Range("E21").Value = "'-2.63%"

$ws.Range("D22").Value = "'0.2418"
$ws.Range("E22").Value = "'3.04%"

$ws.Range("D23").Value = "'0.04370"
$ws.Range("E23").Value = "'-0.63%"

$ws.Range("E24").Value = "'-2.16%"

$ws.Range("D25").Value = "'0.004262"
$ws.Range("E25").Value = "'2.39%"

$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'-0.07%"

$ws.Range("D39").Value = "'0.02032"
$ws.Range("E39").Value = "'-1.26%"

$ws.Range("D40").Value = "'0.05047"
$ws.Range("E40").Value = "'-1.21%"

$ws.Range("D41").Value = "'0.007512"
$ws.Range("E41").Value = "'0.22%"

$ws.Range("D42").Value = "'0.009706"
$ws.Range("E42").Value = "'-4.18%"

$ws.Range("D43").Value = "'0.1341"
$ws.Range("E43").Value = "'-1.88%"

$ws.Range("D44").Value = "'0.002095"
$ws.Range("E44").Value = "'-1.48%"

$ws.Range("D45").Value = "'0.009829"

$ws.Range("D46").Value = "'0.00006198"
$ws.Range("E46").Value = "'-1.44%"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.02%"

$ws.Range("D48").Value = "'0.002877"

$ws.Range("E49").Value = "'12.48%"

$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.02%"

$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.02%"
